# error solve ifrs list
# Corrects the per-period financial metrics on the "company_list" sheet:
# a batch of cells had the wrong (stale) figures, and several metric
# columns that are no longer reported for certain periods are cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1498
$ws.Range("E2").Value = 108
$ws.Range("F2").Value = 108
$ws.Range("G2").Value = 56
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 41
$ws.Range("K2").Value = 969
$ws.Range("L2").Value = 894
$ws.Range("M2").Value = 75
$ws.Range("P2").Value = 62
$ws.Range("Q2").Value = 167
$ws.Range("R2").Value = -64
$ws.Range("S2").Value = -96
$ws.Range("T2").Value = 64
$ws.Range("U2").Value = 103
$ws.Range("V2").Value = 318
$ws.Range("W2").Value = 7.21
$ws.Range("X2").Value = 2.76
$ws.Range("Y2").Value = 96.36
$ws.Range("Z2").Value = 4.5
$ws.Range("AA2").Value = 1188.36
$ws.Range("AB2").Value = -73.08
$ws.Range("AC2").Value = 332
$ws.Range("AE2").Value = 639
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 12449140
$ws.Range("J2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AH2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2093
$ws.Range("E3").Value = 353
$ws.Range("F3").Value = 353
$ws.Range("G3").Value = 336
$ws.Range("H3").Value = 247
$ws.Range("I3").Value = 247
$ws.Range("K3").Value = 1419
$ws.Range("L3").Value = 1089
$ws.Range("M3").Value = 331
$ws.Range("P3").Value = 62
$ws.Range("Q3").Value = 275
$ws.Range("R3").Value = -26
$ws.Range("S3").Value = -90
$ws.Range("T3").Value = 19
$ws.Range("U3").Value = 256
$ws.Range("V3").Value = 252
$ws.Range("W3").Value = 16.88
$ws.Range("X3").Value = 11.8
$ws.Range("Y3").Value = 121.71
$ws.Range("Z3").Value = 20.68
$ws.Range("AA3").Value = 329.28
$ws.Range("AB3").Value = 323.76
$ws.Range("AC3").Value = 1984
$ws.Range("AE3").Value = 2809
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 12449140
$ws.Range("J3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AH3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3776
$ws.Range("E4").Value = 796
$ws.Range("F4").Value = 796
$ws.Range("G4").Value = 819
$ws.Range("H4").Value = 595
$ws.Range("I4").Value = 595
$ws.Range("K4").Value = 2437
$ws.Range("L4").Value = 1524
$ws.Range("M4").Value = 913
$ws.Range("P4").Value = 62
$ws.Range("Q4").Value = 262
$ws.Range("R4").Value = -169
$ws.Range("S4").Value = -25
$ws.Range("T4").Value = 105
$ws.Range("U4").Value = 157
$ws.Range("V4").Value = 248
$ws.Range("W4").Value = 21.07
$ws.Range("X4").Value = 15.76
$ws.Range("Y4").Value = 95.73999999999999
$ws.Range("Z4").Value = 30.87
$ws.Range("AA4").Value = 166.99
$ws.Range("AB4").Value = 1279.9
$ws.Range("AC4").Value = 4781
$ws.Range("AE4").Value = 7752
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 12449140
$ws.Range("J4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("AD4").ClearContents()
$ws.Range("AH4").ClearContents()

# Row 5
$ws.Range("D5").Value = 6028
$ws.Range("E5").Value = 872
$ws.Range("F5").Value = 872
$ws.Range("G5").Value = 657
$ws.Range("H5").Value = 452
$ws.Range("I5").Value = 452
$ws.Range("K5").Value = 3404
$ws.Range("L5").Value = 2088
$ws.Range("M5").Value = 1315
$ws.Range("N5").Value = 1315
$ws.Range("P5").Value = 62
$ws.Range("Q5").Value = 294
$ws.Range("R5").Value = -385
$ws.Range("S5").Value = 272
$ws.Range("T5").Value = 307
$ws.Range("U5").Value = -12
$ws.Range("V5").Value = 798
$ws.Range("W5").Value = 14.47
$ws.Range("X5").Value = 7.49
$ws.Range("Y5").Value = 40.53
$ws.Range("Z5").Value = 15.46
$ws.Range("AA5").Value = 158.76
$ws.Range("AB5").Value = 2193.79
$ws.Range("AC5").Value = 3627
$ws.Range("AE5").Value = 11175
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 12449140
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("AD5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("D6").Value = 6218
$ws.Range("E6").Value = 531
$ws.Range("F6").Value = 531
$ws.Range("G6").Value = 582
$ws.Range("H6").Value = 436
$ws.Range("I6").Value = 436
$ws.Range("K6").Value = 5148
$ws.Range("L6").Value = 3385
$ws.Range("M6").Value = 1762
$ws.Range("N6").Value = 1762
$ws.Range("P6").Value = 62
$ws.Range("Q6").Value = 116
$ws.Range("R6").Value = -887
$ws.Range("S6").Value = 945
$ws.Range("T6").Value = 478
$ws.Range("U6").Value = -363
$ws.Range("V6").Value = 1787
$ws.Range("W6").Value = 8.539999999999999
$ws.Range("X6").Value = 7.01
$ws.Range("Y6").Value = 28.34
$ws.Range("Z6").Value = 10.2
$ws.Range("AA6").Value = 192.08
$ws.Range("AB6").Value = 2885.15
$ws.Range("AC6").Value = 3503
$ws.Range("AE6").Value = 14974
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AJ6").Value = 12449140
$ws.Range("AD6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# Row 7
$ws.Range("D7").Value = 8184
$ws.Range("E7").Value = 1133
$ws.Range("G7").Value = 1100
$ws.Range("H7").Value = 880
$ws.Range("I7").Value = 885
$ws.Range("K7").Value = 6470
$ws.Range("L7").Value = 3830
$ws.Range("M7").Value = 2640
$ws.Range("N7").Value = 2640
$ws.Range("P7").Value = 70
$ws.Range("Q7").Value = 150
$ws.Range("R7").Value = -920
$ws.Range("S7").Value = 260
$ws.Range("U7").Value = 100
$ws.Range("W7").Value = 13.84
$ws.Range("X7").Value = 10.75
$ws.Range("Y7").Value = 40.2
$ws.Range("Z7").Value = 15.15
$ws.Range("AA7").Value = 145.08
$ws.Range("AC7").Value = 6924
$ws.Range("AD7").Value = 15.45
$ws.Range("AE7").Value = 19515
$ws.Range("AF7").Value = 5.48
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("T7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 10058
$ws.Range("E8").Value = 1526
$ws.Range("G8").Value = 1550
$ws.Range("H8").Value = 1210
$ws.Range("I8").Value = 1210
$ws.Range("K8").Value = 7630
$ws.Range("L8").Value = 3790
$ws.Range("M8").Value = 3840
$ws.Range("N8").Value = 3840
$ws.Range("P8").Value = 70
$ws.Range("Q8").Value = 1570
$ws.Range("R8").Value = -1300
$ws.Range("S8").Value = -190
$ws.Range("U8").Value = 100
$ws.Range("W8").Value = 15.17
$ws.Range("X8").Value = 12.03
$ws.Range("Y8").Value = 37.35
$ws.Range("Z8").Value = 17.16
$ws.Range("AA8").Value = 98.7
$ws.Range("AC8").Value = 8517
$ws.Range("AD8").Value = 12.56
$ws.Range("AE8").Value = 28385
$ws.Range("AF8").Value = 3.77
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("T8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 12470
$ws.Range("E9").Value = 1960
$ws.Range("G9").Value = 1990
$ws.Range("H9").Value = 1550
$ws.Range("I9").Value = 1550
$ws.Range("K9").Value = 9060
$ws.Range("L9").Value = 3690
$ws.Range("M9").Value = 5370
$ws.Range("N9").Value = 5370
$ws.Range("P9").Value = 70
$ws.Range("Q9").Value = 1440
$ws.Range("R9").Value = -1300
$ws.Range("S9").Value = -100
$ws.Range("U9").Value = 530
$ws.Range("W9").Value = 15.72
$ws.Range("X9").Value = 12.43
$ws.Range("Y9").Value = 33.66
$ws.Range("Z9").Value = 18.57
$ws.Range("AA9").Value = 68.72
$ws.Range("AC9").Value = 10910
$ws.Range("AD9").Value = 9.81
$ws.Range("AE9").Value = 39695
$ws.Range("AF9").Value = 2.7
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("T9").ClearContents()
$ws.Range("AI9").ClearContents()
